{"js": "const replacements = [\n  [\"14\u00d781=\", \"79\u00d794=\"],\n  [\"92\u00d792=\", \"57\u00d762=\"],\n  [\"64\u00d796=\", \"58\u00d792=\"],\n  [\"37\u00d713=\", \"74\u00d760=\"],\n  [\"43\u00d767=\", \"46\u00d792=\"],\n  [\"93\u00d773=\", \"23\u00d767=\"],\n  [\"32\u00d744=\", \"61\u00d771=\"],\n  [\"96\u00d734=\", \"73\u00d750=\"],\n  [\"38\u00d759=\", \"91\u00d745=\"],\n  [\"35\u00d787=\", \"86\u00d771=\"],\n  [\"97\u00d796=\", \"32\u00d718=\"],\n  [\"81\u00d748=\", \"12\u00d748=\"],\n  [\"63\u00d763=\", \"52\u00d748=\"],\n  [\"41\u00d747=\", \"82\u00d739=\"],\n  [\"38\u00d771=\", \"78\u00d733=\"],\n  [\"50\u00d765=\", \"29\u00d774=\"],\n  [\"89\u00d741=\", \"99\u00d719=\"],\n  [\"70\u00d777=\", \"99\u00d775=\"],\n  [\"13\u00d754=\", \"20\u00d730=\"],\n  [\"27\u00d724=\", \"84\u00d725=\"],\n  [\"56\u00d755=\", \"39\u00d736=\"],\n  [\"78\u00d795=\", \"82\u00d756=\"],\n  [\"38\u00d753=\", \"43\u00d748=\"],\n  [\"69\u00d723=\", \"56\u00d772=\"],\n  [\"22\u00d715=\", \"71\u00d722=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"14\u00d781=\"; New = \"79\u00d794=\" },\n    @{ Old = \"92\u00d792=\"; New = \"57\u00d762=\" },\n    @{ Old = \"64\u00d796=\"; New = \"58\u00d792=\" },\n    @{ Old = \"37\u00d713=\"; New = \"74\u00d760=\" },\n    @{ Old = \"43\u00d767=\"; New = \"46\u00d792=\" },\n    @{ Old = \"93\u00d773=\"; New = \"23\u00d767=\" },\n    @{ Old = \"32\u00d744=\"; New = \"61\u00d771=\" },\n    @{ Old = \"96\u00d734=\"; New = \"73\u00d750=\" },\n    @{ Old = \"38\u00d759=\"; New = \"91\u00d745=\" },\n    @{ Old = \"35\u00d787=\"; New = \"86\u00d771=\" },\n    @{ Old = \"97\u00d796=\"; New = \"32\u00d718=\" },\n    @{ Old = \"81\u00d748=\"; New = \"12\u00d748=\" },\n    @{ Old = \"63\u00d763=\"; New = \"52\u00d748=\" },\n    @{ Old = \"41\u00d747=\"; New = \"82\u00d739=\" },\n    @{ Old = \"38\u00d771=\"; New = \"78\u00d733=\" },\n    @{ Old = \"50\u00d765=\"; New = \"29\u00d774=\" },\n    @{ Old = \"89\u00d741=\"; New = \"99\u00d719=\" },\n    @{ Old = \"70\u00d777=\"; New = \"99\u00d775=\" },\n    @{ Old = \"13\u00d754=\"; New = \"20\u00d730=\" },\n    @{ Old = \"27\u00d724=\"; New = \"84\u00d725=\" },\n    @{ Old = \"56\u00d755=\"; New = \"39\u00d736=\" },\n    @{ Old = \"78\u00d795=\"; New = \"82\u00d756=\" },\n    @{ Old = \"38\u00d753=\"; New = \"43\u00d748=\" },\n    @{ Old = \"69\u00d723=\"; New = \"56\u00d772=\" },\n    @{ Old = \"22\u00d715=\"; New = \"71\u00d722=\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Old\n    $find.Replacement.Text = $r.New\n    $find.Execute([ref]$r.Old, $false, $true, $false, $false, $false, $true, 1, $false, [ref]$r.New, 2)\n}\n"}
